$p = $ppt.ActivePresentation

# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> currently "Integral"      (used by the slide master / slides)
#   ppt/theme/theme2.xml  -> currently "Office Theme"  (used by the notes master)
# The target edit swaps the two themes' contents (their font/format schemes are
# already identical; only the 12 colour-scheme slots actually differ).
# This host's object model only exposes an editable ThemeColorScheme bound to
# theme1.xml (the presentation's active theme), so we push the "Office Theme"
# palette onto it here via ThemeColorScheme.Colors(i).RGB, in the standard
# PowerPoint theme colour order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink

$officeThemeRGB = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
